$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Date
$ws.Range("D3").Value = 44336.358506944445
$ws.Range("E3").Value = 44336.361273148148
$ws.Range("F3").Value = 44336.363518518519

# Row 4 - W
$ws.Range("E4").Value = 32
$ws.Range("F4").Value = 31

# Row 5 - BS
$ws.Range("E5").Value = 1.2
$ws.Range("F5").Value = 1.1

# Row 6 - SE
$ws.Range("E6").Value = 0.1
$ws.Range("F6").Value = 0.9

# Row 7 - SSE
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 1

# Row 8 - Age
$ws.Range("F8").Value = 27

# Row 9
$ws.Range("B9").Value = 3
